# This script updates the "05_abr" sheet data:
#  - removes the "cidade" / "Óbitos" sub-header row (row 102) that separated
#    the confirmed-cases block from the deaths block
#  - removes the "outros estados" / "outros paises" rows (rows 100-101)
#  - removes the "municipios" / "Casos confirmados" sub-header row (row 2)
# so that the deaths data immediately follows the confirmed-cases data,
# with only the single top-level header row ("Unnamed: 0" / "Unnamed: 1") remaining.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Delete rows from the bottom up so earlier row numbers stay valid.
$ws.Rows.Item(102).Delete()
$ws.Rows.Item(101).Delete()
$ws.Rows.Item(100).Delete()
$ws.Rows.Item(2).Delete()
